$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "coupling weight" column (F) values for rows 2-10
$ws.Range("F2").Value = 0.24909999999999999
$ws.Range("F3").Value = 0.20200000000000001
$ws.Range("F4").Value = 0.20200000000000001
$ws.Range("F5").Value = 0.20200000000000001
$ws.Range("F6").Value = 0.19950000000000001
$ws.Range("F7").Value = 0.17599999999999999
$ws.Range("F8").Value = 0.20069999999999999
$ws.Range("F9").Value = 0.2392
$ws.Range("F10").Value = 0.16109999999999999

# Re-enter the "total" column (G) formula across the whole range so Excel
# stores it as a single shared formula (t="shared") like the author's edit
$ws.Range("G2:G10").Formula = "=(A2*E2)-((1-A2)*F2)"

# Update the active selection to match the saved cursor position
$ws.Range("F11").Select()
